# Update gh-pages to output generated at 456a3b4
# Updates the "想去人数" (want-to-go count) figures in column F for the
# "展览" (Worksheets index 1) and "全部类型" (Worksheets index 4) sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (1st sheet) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value = 372
$ws1.Range("F5").Value = 419
$ws1.Range("F6").Value = 265
$ws1.Range("F7").Value = 2413
$ws1.Range("F8").Value = 417
$ws1.Range("F9").Value = 6298

# --- Sheet "全部类型" (4th sheet) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 372
$ws4.Range("F5").Value = 419
$ws4.Range("F6").Value = 265
$ws4.Range("F9").Value = 2413
$ws4.Range("F10").Value = 417
$ws4.Range("F11").Value = 6298
